$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("T", "V", "X")
for ($row = 2; $row -le 6; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $text = $cell.Value()
        if ($text -and $text.Length -gt 0) {
            $cell.Value = $text.Substring(0,1).ToUpper() + $text.Substring(1)
        }
    }
}

$ws.Range("S1:S6").NumberFormat = "#,##0"
$ws.Range("U1:U6").NumberFormat = "#,##0"
$ws.Range("W1:W6").NumberFormat = "#,##0"
